# Update the Agenda table on slide 2: add "Sammy Douglas" row and resize the
# table/fonts so the (now 13-row) table still fits the available space.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the Agenda table shape (named "Table 5" in the source deck).
$tblShape = $null
foreach ($shp in $s.Shapes) {
    if ($shp.HasTable) {
        $tblShape = $shp
    }
}

$tbl = $tblShape.Table

# Insert a new row right after "Anita Wong" (row 7), before "Senthil/Melissa".
$newRow = $tbl.Rows.Add(8)
$tbl.Cell(8, 1).Shape.TextFrame.TextRange.Text = "Sammy Douglas"
$tbl.Cell(8, 2).Shape.TextFrame.TextRange.Text = "Destruction Update"
$tbl.Cell(8, 3).Shape.TextFrame.TextRange.Text = "5 minutes"

# Shrink the header font (32 -> 28 pt) and body font (28 -> 24 pt) so the
# extra row still fits, matching the rest of the Agenda rows.
for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
    $tbl.Cell(1, $c).Shape.TextFrame.TextRange.Font.Size = 28
}
for ($r = 2; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $tbl.Cell($r, $c).Shape.TextFrame.TextRange.Font.Size = 24
    }
}

# Resize/reposition the table frame and columns to their final layout
# *before* tightening row heights, since changing the frame size causes
# PowerPoint to redistribute row heights proportionally.
$tblShape.Left = 13.57755905511811
$tblShape.Top = 110.53385826771654
$tblShape.Width = 932.8448031496063
$tblShape.Height = 416.8255905511811

$tbl.Columns.Item(1).Width = 203.70275590551182
$tbl.Columns.Item(2).Width = 563.5144094488189
$tbl.Columns.Item(3).Width = 165.6276377952756

# Tighten row heights to compensate for the newly inserted row.
$tbl.Rows.Item(1).Height = 33.32779527559055
for ($r = 2; $r -le ($tbl.Rows.Count - 1); $r++) {
    $tbl.Rows.Item($r).Height = 31.740708661417322
}
$tbl.Rows.Item($tbl.Rows.Count).Height = 33.32779527559055

$tblShape.Name = "Table 2"
